$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan")
$ws.Activate()

# The "Data Structures and Algorithms" course (CSCA 5424, row 38 - planned for
# 2026 Spring 1) was completed, so it moves up into the "2025 Fall 2" block
# (row 35), swapping places with "Machine Learning" (DTSA 5511), which moves
# out to the later "2026 Spring 1" term (row 38). Swap the two rows' data,
# including the course-color fill on column A, and let the SUM formulas in
# F32/F36 recalculate automatically.

$a35 = $ws.Range("A35").Value2
$b35 = $ws.Range("B35").Value2
$c35 = $ws.Range("C35").Value2
$d35 = $ws.Range("D35").Value2
$e35 = $ws.Range("E35").Value2
$f35 = $ws.Range("F35").Value2
$color35 = $ws.Range("A35").Interior.Color

$a38 = $ws.Range("A38").Value2
$b38 = $ws.Range("B38").Value2
$c38 = $ws.Range("C38").Value2
$d38 = $ws.Range("D38").Value2
$e38 = $ws.Range("E38").Value2
$f38 = $ws.Range("F38").Value2
$color38 = $ws.Range("A38").Interior.Color

$ws.Range("A35").Value2 = $a38
$ws.Range("B35").Value2 = $b38
$ws.Range("C35").Value2 = $c38
$ws.Range("D35").Value2 = $d38
$ws.Range("E35").Value2 = $e38
$ws.Range("F35").Value2 = $f38
$ws.Range("A35").Interior.Color = $color38

$ws.Range("A38").Value2 = $a35
$ws.Range("B38").Value2 = $b35
$ws.Range("C38").Value2 = $c35
$ws.Range("D38").Value2 = $d35
$ws.Range("E38").Value2 = $e35
$ws.Range("F38").Value2 = $f35
$ws.Range("A38").Interior.Color = $color35

# Update the view to match where the user ended up scrolled/selected.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F33").Select()
